$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-09-17 (45186) to 2023-09-19 (45188) for every data row (rows 2-138).
for ($row = 2; $row -le 138; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
